$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "TestConfig" row label to "A"
$ws.Range("A2").Value = "A"

# Add a new config row (row 4) labeled "B" with its parameter values
$ws.Range("A4").Value = "B"
$ws.Range("B4").Value = 150
$ws.Range("C4").Value = 600
$ws.Range("D4").Value = 300
$ws.Range("E4").Value = 300
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 500

# F2/F3 use a quote-prefixed number format; copy that formatting onto F4
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the selection to J4 to mirror the saved cursor position
$ws.Range("J4").Select()
